# Updated integration test spreadsheets with content_description fields.
# Fixes #1033.
#
# This adds three new columns (content_description.text / .ontology /
# .ontology_label) to the "Sequence file" sheet, between the existing
# "INSDC run accession" column (H) and the "Protocol ID" columns
# (old I/J, which shift right to L/M), and makes "Sequence file" the
# active/selected sheet (it had been "Project - Contributors").

$wb = $excel.ActiveWorkbook

$wsSeq = $wb.Worksheets.Item("Sequence file")
$wsContrib = $wb.Worksheets.Item("Project - Contributors")

# The old J5:K5 merged (empty) cell needs to end up unmerged after the
# insert, with the old style widened over two plain cells. Unmerge first
# so the shift carries a clean, unmerged pair of cells into the new
# L5:M5 area under the column insert below -- no, these shift to M5:N5.
$wsSeq.Range("J5:K5").UnMerge()

# Insert three blank columns at I:K. Excel shifts the old I/J/K (Protocol
# ID / Protocol ID / Process ID columns) right to L/M/N, carrying their
# styles and values with them, and the new I/J/K inherit the row
# formatting to their left.
$wsSeq.Columns("I:K").Insert()

# --- Row 1 (field descriptions) ---
$wsSeq.Cells.Item(1, 9).Value = "General description of the contents of the file."
$wsSeq.Cells.Item(1, 10).Value = "An ontology term identifier in the form prefix:accession."
$wsSeq.Cells.Item(1, 11).Value = "The preferred label for the ontology term referred to in the ontology field. This may differ from the user-supplied value in the text field."

# --- Row 2 (short field labels) ---
$wsSeq.Cells.Item(2, 9).Value = "Content description"
$wsSeq.Cells.Item(2, 10).Value = "Content description ontology"
$wsSeq.Cells.Item(2, 11).Value = "Content description ontology label"

# --- Row 3 (example values) ---
$wsSeq.Cells.Item(3, 9).Value = "DNA sequence (raw); Sequence alignment"
$wsSeq.Cells.Item(3, 10).Value = "DATA:3497; DATA:0863"
$wsSeq.Cells.Item(3, 11).Value = "DNA sequence (raw); Sequence alignment"

# --- Row 4 (programmatic field tags) ---
$wsSeq.Cells.Item(4, 9).Value = "sequence_file.file_core.content_description.text"
$wsSeq.Cells.Item(4, 10).Value = "sequence_file.file_core.content_description.ontology"
$wsSeq.Cells.Item(4, 11).Value = "sequence_file.file_core.content_description.ontology_label"

# --- Rows 6-8 (sample data rows) ---
foreach ($r in 6..8) {
    $wsSeq.Cells.Item($r, 9).Value = "DNA sequence (raw)"
    $wsSeq.Cells.Item($r, 10).Value = "data:3497"
    $wsSeq.Cells.Item($r, 11).Value = "DNA sequence (raw)"
}

# Make "Sequence file" the selected/active sheet (was "Project - Contributors").
$wsSeq.Select()
$wsSeq.Range("I6:K6").Select()
